$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------------
# 1) Remove the two rows that disappear entirely ("Bomb & Brick FX(Optional)"
#    and "Sounds(optional)" - old rows 18 & 19). This shifts the old Total
#    row (20) up to row 18 and automatically repairs the SUM() ranges.
# ---------------------------------------------------------------------------
$ws.Rows("18:19").Delete()

# ---------------------------------------------------------------------------
# 2) Grab a "border-only" format sample (old style used by I11) and stamp it
#    onto the new column F filler cells (rows 5..18, now that the Total row
#    has settled at row 18).
# ---------------------------------------------------------------------------
$ws.Range("I11").Copy()
$ws.Range("F5:F18").PasteSpecial(-4122)

# ---------------------------------------------------------------------------
# 3) Unmerge H10:I10 and drop the old "Legends" label plus its neighbour
#    cell entirely - the new layout has nothing in column H.
# ---------------------------------------------------------------------------
$ws.Range("H10:I10").UnMerge()
$ws.Range("H10").Clear()
$ws.Range("H11").Clear()

# ---------------------------------------------------------------------------
# 4) Re-purpose I9..I12 as the "Remarks" side list, matching the border +
#    centred formatting already used on the task column (C10 is a ready
#    template for that look).
# ---------------------------------------------------------------------------
$ws.Range("C10").Copy()
$ws.Range("I9").PasteSpecial(-4122)
$ws.Range("I9").Value = "AI which works as Player"

$ws.Range("C10").Copy()
$ws.Range("I10").PasteSpecial(-4122)
$ws.Range("I10").Value = "Neutral characters"

$ws.Range("C10").Copy()
$ws.Range("I11").PasteSpecial(-4122)
$ws.Range("I11").Value = "Sounds"

$ws.Range("C10").Copy()
$ws.Range("I12").PasteSpecial(-4122)
$ws.Range("I12").Value = "Effects"

# The "Future Task" callout sits in the header-style (bold) formatting, like
# the column headers in row 4.
$ws.Range("C4").Copy()
$ws.Range("I8").PasteSpecial(-4122)
$ws.Range("I8").Value = "Future Task"

# ---------------------------------------------------------------------------
# 5) New "Remarks" column (F). Header gets the bold header look, the body
#    cells (already stamped with the bordered look in step 1) stay blank.
# ---------------------------------------------------------------------------
$ws.Range("C4").Copy()
$ws.Range("F4").PasteSpecial(-4122)
$ws.Range("F4").Value = "Remarks"

# ---------------------------------------------------------------------------
# 6) Fill in the previously-blank Actual Time values for the last few tasks.
# ---------------------------------------------------------------------------
$ws.Range("E14").Value = 0.5
$ws.Range("E15").Value = 1
$ws.Range("E16").Value = 1.5

# ---------------------------------------------------------------------------
# 7) Column widths: new Remarks column, and the widened callout column I.
# ---------------------------------------------------------------------------
$ws.Columns("F").ColumnWidth = 13.75
$ws.Columns("I").ColumnWidth = 24.25

# ---------------------------------------------------------------------------
# 8) Row heights: row 8 grows to match the header-style entry that now lives
#    in I8; row 12 keeps its wrapped height but becomes an explicit height.
# ---------------------------------------------------------------------------
$ws.Rows(8).RowHeight = 18.75
$ws.Rows(12).RowHeight = 30

# ---------------------------------------------------------------------------
# 9) Selection / view - land on E17:F17 with no frozen scroll offset.
# ---------------------------------------------------------------------------
$ws.Range("E17:F17").Select()

$wb.Application.CalculateFull()
